$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.062.18"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "1.797.99"
$ws.Range("E3").Value = "  -2.46%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "'307.63"
$ws.Range("E5").Value = "  -2.50%  "
$ws.Range("D7").Value = "'0.4207"
$ws.Range("E7").Value = "  -2.35%  "
$ws.Range("D8").Value = "'0.3593"
$ws.Range("E8").Value = "  -2.55%  "
$ws.Range("D9").Value = "'0.07116"
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("D10").Value = "'0.8435"
$ws.Range("E10").Value = "  -3.52%  "
$ws.Range("E11").Value = "  -3.68%  "
$ws.Range("D12").Value = "1.825.41"
$ws.Range("E12").Value = "  -2.88%  "
$ws.Range("E13").Value = "  -3.30%  "
$ws.Range("D14").Value = "'6.363"
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "'0.06765"
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("D16").Value = "'1.006"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "'80.08"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").Value = "'0.000008730"
$ws.Range("E18").Value = "  -3.77%  "
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "'15.01"
$ws.Range("E20").Value = "  -3.54%  "
$ws.Range("D21").Value = "27.073.41"
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D22").Value = "'5.058"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").Value = "'10.97"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "2.026.50"
$ws.Range("E24").Value = "  -3.29%  "
$ws.Range("D25").Value = "'1.927"
$ws.Range("E25").Value = "  -3.04%  "
$ws.Range("D26").Value = "'152.92"
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("D27").Value = "'18.11"
$ws.Range("E27").Value = "  -4.69%  "
$ws.Range("D28").Value = "'5.015"
$ws.Range("E28").Value = "  -5.91%  "
$ws.Range("D29").Value = "'113.38"
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("D30").Value = "'1.644"
$ws.Range("E30").Value = "  -12.51%  "
$ws.Range("D31").Value = "'0.08983"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").Value = "'0.7233"
$ws.Range("E32").Value = "  -7.82%  "
$ws.Range("D33").Value = "'2.869"
$ws.Range("E33").Value = "  -3.81%  "
$ws.Range("D34").Value = "'4.322"
$ws.Range("E34").Value = "  -6.20%  "
$ws.Range("D35").Value = "'1.092"
$ws.Range("E35").Value = "  -6.26%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "'1.078"
$ws.Range("E37").Value = "  -2.88%  "
$ws.Range("D38").Value = "'0.01904"
$ws.Range("E38").Value = "  -3.01%  "
$ws.Range("D39").Value = "'0.05137"
$ws.Range("E39").Value = "  -5.50%  "
$ws.Range("E40").Value = "  -3.81%  "
$ws.Range("D41").Value = "'0.4968"
$ws.Range("E41").Value = "  -3.98%  "
$ws.Range("D42").Value = "'2.605"
$ws.Range("E42").Value = "  -7.89%  "
$ws.Range("D43").Value = "'8.047"
$ws.Range("E43").Value = "  -6.74%  "
$ws.Range("D44").Value = "'5.905"
$ws.Range("E44").Value = "  -12.65%  "
$ws.Range("D45").Value = "'105.14"
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.16"
$ws.Range("E47").Value = "  -4.77%  "
$ws.Range("D48").Value = "'0.06293"
$ws.Range("E48").Value = "  -3.86%  "
$ws.Range("D49").Value = "'0.4513"
$ws.Range("E49").Value = "  -5.82%  "
$ws.Range("D50").Value = "'1.600"
$ws.Range("E50").Value = "  -3.80%  "
$ws.Range("D51").Value = "'1.716"
$ws.Range("E51").Value = "  -7.40%  "
